# Update Name of Algo
# Updates a batch of imputed numeric values in the result_data_RandomForest
# sheet to reflect the latest algorithm run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 5.867900000000002
$ws.Range("D6").Value = -8.059100000000003
$ws.Range("B14").Value = 5.293000000000002
$ws.Range("D18").Value = -8.535300000000001
$ws.Range("D19").Value = -8.399899999999995
$ws.Range("B21").Value = 8.739599999999999
$ws.Range("C22").Value = -12.04410000000001
$ws.Range("B23").Value = 9.208099999999996
$ws.Range("C24").Value = -13.60209999999999
$ws.Range("B25").Value = 5.557000000000001
$ws.Range("B26").Value = 4.766200000000002
$ws.Range("C28").Value = -13.944
$ws.Range("B29").Value = 5.224500000000003
$ws.Range("C36").Value = -11.60200000000001
$ws.Range("D44").Value = -6.340800000000002
$ws.Range("C45").Value = -13.84189999999999
$ws.Range("D47").Value = -7.467600000000001
$ws.Range("C48").Value = -11.2454
$ws.Range("C49").Value = -13.75249999999999
$ws.Range("D51").Value = -8.610499999999996
$ws.Range("C52").Value = -10.8014
$ws.Range("B53").Value = 5.4779
$ws.Range("C53").Value = -10.9094
$ws.Range("C54").Value = -13.66349999999999
$ws.Range("D55").Value = -8.686100000000001
$ws.Range("B57").Value = 4.546299999999996
$ws.Range("D57").Value = -8.344299999999997
$ws.Range("B59").Value = 4.659399999999996
$ws.Range("D64").Value = -8.178599999999987
$ws.Range("B69").Value = 5.672699999999996
$ws.Range("C70").Value = -11.3343
$ws.Range("B79").Value = 9.429600000000006
$ws.Range("D80").Value = -7.991500000000002
$ws.Range("B83").Value = 5.762399999999999
$ws.Range("C86").Value = -13.515
$ws.Range("C87").Value = -12.7003
$ws.Range("C89").Value = -13.32569999999999
$ws.Range("B91").Value = 5.761600000000003
$ws.Range("D92").Value = -7.161000000000004
$ws.Range("B93").Value = 8.605900000000002
$ws.Range("D94").Value = -6.436900000000004
$ws.Range("D96").Value = -8.559799999999997
$ws.Range("C101").Value = -12.3561
$ws.Range("D101").Value = -8.226000000000004
$ws.Range("B103").Value = 5.653400000000008
